$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.337.32"
Set-TextValue $ws.Range("E2") "  +1.92%  "
Set-TextValue $ws.Range("D3") "3.344.48"
Set-TextValue $ws.Range("E3") "  +3.13%  "
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "192.08"
Set-TextValue $ws.Range("E5") "  +3.97%  "
Set-TextValue $ws.Range("D6") "592.70"
Set-TextValue $ws.Range("E6") "  +1.95%  "
Set-TextValue $ws.Range("E7") "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.607"
Set-TextValue $ws.Range("E8") "  +0.86%  "
Set-TextValue $ws.Range("D9") "0.133"
Set-TextValue $ws.Range("E9") "  +2.65%  "
Set-TextValue $ws.Range("D10") "6.71"
Set-TextValue $ws.Range("E10") "  +1.47%  "
Set-TextValue $ws.Range("D11") "0.423"
Set-TextValue $ws.Range("E11") "  +1.67%  "
Set-TextValue $ws.Range("D12") "3.925.75"
Set-TextValue $ws.Range("E12") "  +3.03%  "
Set-TextValue $ws.Range("E13") "  +0.81%  "
Set-TextValue $ws.Range("D14") "28.20"
Set-TextValue $ws.Range("E14") "  +1.41%  "
Set-TextValue $ws.Range("D15") "69.368.38"
Set-TextValue $ws.Range("E15") "  +1.95%  "
Set-TextValue $ws.Range("D16") "0.0000171"
Set-TextValue $ws.Range("E16") "  +0.97%  "
Set-TextValue $ws.Range("D17") "3.339.43"
Set-TextValue $ws.Range("E17") "  +2.93%  "
Set-TextValue $ws.Range("D18") "5.82"
Set-TextValue $ws.Range("E18") "  +0.37%  "
Set-TextValue $ws.Range("D19") "13.70"
Set-TextValue $ws.Range("E19") "  +1.74%  "
Set-TextValue $ws.Range("D20") "425.55"
Set-TextValue $ws.Range("E20") "  +7.30%  "
Set-TextValue $ws.Range("E21") "  +1.76%  "
Set-TextValue $ws.Range("D22") "73.35"
Set-TextValue $ws.Range("E22") "  +2.80%  "
Set-TextValue $ws.Range("E23") "  +0.00%  "
Set-TextValue $ws.Range("D24") "0.517"
Set-TextValue $ws.Range("E24") "  +0.33%  "
Set-TextValue $ws.Range("E25") "  +2.25%  "
Set-TextValue $ws.Range("E26") "  +2.30%  "
Set-TextValue $ws.Range("D27") "9.59"
Set-TextValue $ws.Range("E27") "  -0.34%  "
Set-TextValue $ws.Range("D28") "1.00"
Set-TextValue $ws.Range("E28") "  +0.08%  "
Set-TextValue $ws.Range("E29") "  +2.54%  "
Set-TextValue $ws.Range("D30") "5.62"
Set-TextValue $ws.Range("E30") "  +0.38%  "
Set-TextValue $ws.Range("D31") "22.98"
Set-TextValue $ws.Range("E31") "  +0.85%  "
Set-TextValue $ws.Range("E32") "  +1.20%  "
Set-TextValue $ws.Range("D33") "7.00"
Set-TextValue $ws.Range("E33") "  -0.26%  "
Set-TextValue $ws.Range("D35") "164.92"
Set-TextValue $ws.Range("E35") "  +1.89%  "
Set-TextValue $ws.Range("E36") "  +1.27%  "
Set-TextValue $ws.Range("E37") "  +1.03%  "
Set-TextValue $ws.Range("D38") "27.00"
Set-TextValue $ws.Range("E38") "  +1.38%  "
Set-TextValue $ws.Range("D39") "0.808"
Set-TextValue $ws.Range("E39") "  -0.42%  "
Set-TextValue $ws.Range("D40") "4.57"
Set-TextValue $ws.Range("E40") "  -0.38%  "
Set-TextValue $ws.Range("D41") "2.747.44"
Set-TextValue $ws.Range("E41") "  +5.26%  "
Set-TextValue $ws.Range("D42") "6.45"
Set-TextValue $ws.Range("E42") "  -0.49%  "
Set-TextValue $ws.Range("D43") "2.52"
Set-TextValue $ws.Range("E43") "  +1.33%  "
Set-TextValue $ws.Range("E44") "  +1.33%  "
Set-TextValue $ws.Range("D45") "41.09"
Set-TextValue $ws.Range("E45") "  -0.29%  "
Set-TextValue $ws.Range("D46") "0.0685"
Set-TextValue $ws.Range("E46") "  +0.10%  "
Set-TextValue $ws.Range("D47") "343.36"
Set-TextValue $ws.Range("E47") "  +2.13%  "
Set-TextValue $ws.Range("E48") "  +1.08%  "
Set-TextValue $ws.Range("D49") "32.42"
Set-TextValue $ws.Range("E49") "  +4.20%  "
Set-TextValue $ws.Range("D50") "1.00"
Set-TextValue $ws.Range("E50") "  +2.79%  "
Set-TextValue $ws.Range("D51") "6.28"
Set-TextValue $ws.Range("E51") "  -0.50%  "
